$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44314
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = 13000
$ws.Range("O2").Value = 13000
$ws.Range("P2").Value = 13000
$ws.Range("S2").Value = 1857

# Row 3
$ws.Range("D3").Value = 44314
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 45
$ws.Range("N3").Value = 11000
$ws.Range("O3").Value = 11000
$ws.Range("P3").Value = 11000
$ws.Range("S3").Value = 1571

# Row 4
$ws.Range("D4").Value = 44315
$ws.Range("L4").Value = "Especial"
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 14000
$ws.Range("P4").Value = 14000
$ws.Range("S4").Value = 2000

# Row 5
$ws.Range("D5").Value = 44315
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = 12000
$ws.Range("P5").Value = 12500
$ws.Range("S5").Value = 1786

# Row 6
$ws.Range("D6").Value = 44315
$ws.Range("M6").Value = 80
$ws.Range("N6").Value = 10000
$ws.Range("P6").Value = 10500
$ws.Range("S6").Value = 1500

# Row 7
$ws.Range("D7").Value = 44344

# Row 8
$ws.Range("D8").Value = 44321
$ws.Range("M8").Value = 140
$ws.Range("N8").Value = 11000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 11500
$ws.Range("R8").Value = "Región Metropolitana"
$ws.Range("S8").Value = 1643

# Row 9
$ws.Range("D9").Value = 44321
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 80
$ws.Range("N9").Value = 8000
$ws.Range("O9").Value = 8000
$ws.Range("P9").Value = 8000
$ws.Range("S9").Value = 1143

# Row 10
$ws.Range("D10").Value = 44342
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 50
$ws.Range("O10").Value = 12000
$ws.Range("P10").Value = 12000
$ws.Range("S10").Value = 1714

# Row 11
$ws.Range("D11").Value = 44307
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 70
$ws.Range("N11").Value = 14000
$ws.Range("O11").Value = 14000
$ws.Range("P11").Value = 14000
$ws.Range("S11").Value = 2000

# Row 12
$ws.Range("D12").Value = 44307
$ws.Range("L12").Value = "Segunda"
$ws.Range("N12").Value = 10000
$ws.Range("O12").Value = 10000
$ws.Range("P12").Value = 10000
$ws.Range("S12").Value = 1429

# Row 13
$ws.Range("D13").Value = 44302
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 340
$ws.Range("N13").Value = 12000
$ws.Range("O13").Value = 13000
$ws.Range("P13").Value = 12500
$ws.Range("R13").Value = "Provincia de Santiago"
$ws.Range("S13").Value = 1786

# Row 15
$ws.Range("D15").Value = 44300
$ws.Range("M15").Value = 150
$ws.Range("N15").Value = 12000
$ws.Range("P15").Value = 12500
$ws.Range("R15").Value = "Provincia de Santiago"
$ws.Range("S15").Value = 1786

# Row 16
$ws.Range("D16").Value = 44335
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 80
$ws.Range("N16").Value = 14000
$ws.Range("O16").Value = 14000
$ws.Range("P16").Value = 14000
$ws.Range("S16").Value = 2000

# Row 19
$ws.Range("D19").Value = 44316
$ws.Range("M19").Value = 40
$ws.Range("N19").Value = 13000
$ws.Range("O19").Value = 13000
$ws.Range("P19").Value = 13000
$ws.Range("S19").Value = 1857

# Row 20
$ws.Range("D20").Value = 44316
$ws.Range("L20").Value = "Segunda"
$ws.Range("M20").Value = 50
$ws.Range("O20").Value = 11000
$ws.Range("P20").Value = 11000
$ws.Range("S20").Value = 1571

# Row 21
$ws.Range("D21").Value = 44322
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 11000
$ws.Range("O21").Value = 11000
$ws.Range("P21").Value = 11000
$ws.Range("S21").Value = 1571

# Row 22
$ws.Range("D22").Value = 44312
$ws.Range("M22").Value = 50
$ws.Range("N22").Value = 13000
$ws.Range("P22").Value = 13000
$ws.Range("R22").Value = "Región Metropolitana"
$ws.Range("S22").Value = 1857

# Row 23
$ws.Range("D23").Value = 44312
$ws.Range("M23").Value = 20
$ws.Range("N23").Value = 11000
$ws.Range("O23").Value = 11000
$ws.Range("P23").Value = 11000
$ws.Range("S23").Value = 1571
